# Update column G ("K" - strikeouts) with regenerated values, replacing the
# old "Strike#" derived values with the new K-based values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$newK = @{
    2  = 2
    3  = 2
    4  = 5
    5  = 6
    6  = 4
    7  = 5
    8  = 4
    9  = 6
    10 = 10
    11 = 5
    12 = 4
    13 = 3
    14 = 4
    15 = 9
    16 = 4
    17 = 6
    18 = 5
    19 = 0
    20 = 3
    21 = 6
    22 = 4
    23 = 3
    24 = 9
    25 = 5
    26 = 8
    27 = 8
    28 = 5
    29 = 4
    30 = 9
    31 = 5
    32 = 4
    33 = 4
    34 = 1
}

foreach ($row in $newK.Keys) {
    $ws.Cells.Item($row, 7).Value = $newK[$row]
}
